$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare formatting for the 18 new rows (25:42) by copying the
# existing column-A / column-B formats (A column = bold/centered/bordered
# style; B column = plain default style) before writing any values. ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A25:A42").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B25:B42").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Write / overwrite the primer id + sequence values for rows 1:42. ---
$values = @(
    ("dxs_sub_933_C_T_sub-1", "GTAAGCTCCTGGGGCTTCAC"),
    ("dxs_sub_933_C_T_sub-2", "GGCCTCGTGATACGCCTATTACCGAAAAGTAGCGGCGGTTT"),
    ("dxs_sub_933_C_T_sub-3", "CAAATAAAACGAAAGGCTCATAGCTCGGCAAACCGCCGCTACTTTTCGGTAAACAACCGCTGGAGGGATC"),
    ("dxs_sub_933_C_T_sub-4", "CAGTCACTCGATACCTCGGC"),
    ("dxs_sub_933_C_T_sub-7", "AATAGGCGTATCACGAGGCC"),
    ("dxs_sub_933_C_T_sub-8", "TGAGCCTTTCGTTTTATTTGATGC"),
    ("test-dxs_sub_933_C_T_sub-1", "AACGCGAAGGTCGGTTTTTC"),
    ("test-dxs_sub_933_C_T_sub-2", "TGCCGACCAGCAACTTGG"),
    ("tktA_573_12bp_sub-1", "AGGTTCTGACGGGAGAGGAT"),
    ("tktA_573_12bp_sub-2", "GGCCTCGTGATACGCCTATTCATGTCGAGGGATGGTTCACCGACGACACC"),
    ("tktA_573_12bp_sub-3", "CAAATAAAACGAAAGGCTCAGGTGTCGTCGGTGAACCATCCCTCGACATGACCATCGATAGAAATACCGTTGTCA"),
    ("tktA_573_12bp_sub-4", "CACCTGCGCATCACTCAAAC"),
    ("tktA_573_12bp_sub-7", "AATAGGCGTATCACGAGGCC"),
    ("tktA_573_12bp_sub-8", "TGAGCCTTTCGTTTTATTTGATGC"),
    ("test-tktA_573_12bp_sub-1", "AACTTCTGAACCGGTAGCGA"),
    ("test-tktA_573_12bp_sub-2", "AGCGCATTGGCAATAGTTGTC"),
    ("aceE_del-1", "GCGTCACAGACATGAAATTGGT"),
    ("aceE_del-2", "GGCCTCGTGATACGCCTATTGGGTTATTCCTTATCTATCT"),
    ("aceE_del-3", "CAAATAAAACGAAAGGCTCACAACGTTATTAGATAGATAAGGAATAACCCGAGGTAAAAGAATAATGGCTATCGA"),
    ("aceE_del-4", "TTAACACCAAACTCGCGTGC"),
    ("aceE_del-7", "AATAGGCGTATCACGAGGCC"),
    ("aceE_del-8", "TGAGCCTTTCGTTTTATTTGATGC"),
    ("test-aceE_del-1", "ACGTAAAGTCTACATTTGTGCA"),
    ("test-aceE_del-2", "GGAGCTGCTTCTGCACGTTT"),
    ("pntA_promoter_sub-1", "CGAGGTTTGTGCCGTAAAGC"),
    ("pntA_promoter_sub-2", "GGCCTCGTGATACGCCTATTGTCCTAGGTATAATGCTAGCACGAATCTAGAGAAAGATTGGACGTACCATAATGCGAATTGGCATACCAAGAG"),
    ("pntA_promoter_sub-3", "CAAATAAAACGAAAGGCTCACTAGATTCGTGCTAGCATTATACCTAGGACTGAGCTAGCTGTCAAGGCGCGGTGATAGTGGGATAAACACCT"),
    ("pntA_promoter_sub-4", "ACTTGGTGATGCGGTAGTCG"),
    ("pntA_promoter_sub-7", "AATAGGCGTATCACGAGGCC"),
    ("pntA_promoter_sub-8", "TGAGCCTTTCGTTTTATTTGATGC"),
    ("test-pntA_promoter_sub-1", "TAATTTCGCCCGCACGGAT"),
    ("test-pntA_promoter_sub-2", "TGTCGAACGGGACCATCATC"),
    ("Cgl1452_ins-1", "CACTGCGCGGGATTTTATGG"),
    ("Cgl1452_ins-2", "TCAATACTCTTTTTGGCGCGCATGTGAACGCCTGACCAGG"),
    ("Cgl1452_ins-3", "CAAATAAAACGAAAGGCTCATCCGGCGACCGCTCCGAGGTTGAAGCTTAAGCATCCGGCATGAACAAAGC"),
    ("Cgl1452_ins-4", "CGATGTCGCTGGCGTTAATG"),
    ("Cgl1452_ins-7", "AATAGGCGTATCACGAGGCC"),
    ("Cgl1452_ins-8", "TGAGCCTTTCGTTTTATTTGATGC"),
    ("Cgl1452_ins-5", "CGCGCCAAAAAGAGTATTGACT"),
    ("Cgl1452_ins-6", "GGCCTCGTGATACGCCTATTTTAAGCTTCAACCTCGGAGCG"),
    ("test-Cgl1452_ins-1", "AGTCGCTAAAGTCAGGCCAT"),
    ("test-Cgl1452_ins-2", "TGACTTGTTAGCCGGTCAGC")
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
